$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A ("No.") is removed entirely; every other column shifts
# one position to the left (B->A, C->B, D->C, E->D, F->E).
$ws.Range("A1").EntireColumn.Delete() | Out-Null

# New trailing column F holds a "Video" header, styled to match the other
# bold header cells (e.g. copy E1's bold formatting across).
$ws.Range("E1").Copy($ws.Range("F1")) | Out-Null
$ws.Range("F1").Value = "Video"

# Column width tweaks (post shift).
$ws.Columns("A").ColumnWidth = 15.5
$ws.Columns("B").ColumnWidth = 15.833333333333334
$ws.Columns("E").ColumnWidth = 10.333333333333334

# Selection / view bookkeeping.
$ws.Range("E5").Select() | Out-Null

# Window was maximized/repositioned on a different monitor since the file
# was last saved.
$win = $wb.Windows.Item(1)
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12456
